# Fix the "Requirement ID" numbering in the Requirements sheet (column A).
# A number of rows (roughly rows 25-269) had stale/duplicated sequence
# numbers in column A; renumber them so each row's value is simply
# (row number - 1), matching the surrounding, already-correct rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

for ($r = 25; $r -le 269; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $expected = $r - 1
    if ($cell.Value -ne $expected) {
        $cell.Value = $expected
    }
}

# Restore the final on-screen viewport/selection state that was left behind
# after the fix-up pass: frozen-pane top-left cell and the highlighted
# range of the last rows that were corrected.
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 394
$ws.Range("A265:A269").Select()
